# Update Streamlit app with latest changes
# Append 4 new rows of EDM DATA sensor readings (rows 200-203).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2024-12-19 00:28:04", -0.120646696924232,  -0.001781832359011997, 0.008598887543500385),
    @("2024-12-19 00:28:05", -0.1202628043195,     -0.001749642003835996, 0.008416674157460247),
    @("2024-12-19 00:28:06", -0.1203181108811986,  -0.001944100092807996, 0.009356418021224843),
    @("2024-12-19 00:28:07", -0.1209069630969316,  -0.001833812696693996, 0.008868828961834625)
)

$startRow = 200
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
}
